$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 헌병/운전병 short labels (column E) into the fuller descriptions
$ws.Range("E12").Value = "군사경찰(헌병)"
$ws.Range("E13").Value = "운전병,크레인차량운전"

# Remove the stray leftover cells (C17/D17 "설명"/"사진", B18 "소총수")
$ws.Range("C17").ClearContents()
$ws.Range("D17").ClearContents()
$ws.Range("B18").ClearContents()

# Add three new duty-type rows to the B/D/E reference table
$ws.Range("B14").Value = "인사"
$ws.Range("D14").Value = "인사"
$ws.Range("E14").Value = "일반행정병"

$ws.Range("B15").Value = "군종"
$ws.Range("D15").Value = "군종"
$ws.Range("E6").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = "군종병"

$ws.Range("B16").Value = "기타"
$ws.Range("D16").Value = "기타"
$ws.Range("E16").Value = "조리병, "

# Append the new survey questions below the existing question list
$ws.Range("B44").Value = "혼밥을 얼마나 자주 하십니까"
$ws.Range("B45").Value = "운동을 주3회이상 한다."
$ws.Range("B46").Value = "눈치가 빠르다 느리다"
$ws.Range("B47").Value = "비밀을 잘 지킨다."
$ws.Range("B48").Value = "참을성"
$ws.Range("B49").Value = "법을 잘지킨다"
$ws.Range("B50").Value = "맛없는 반찬 맛있게먹는다"
$ws.Range("B51").Value = "나는 된장찌개를 끓일줄 안다"
$ws.Range("B52").Value = "맛소금과 소금의 차이를 안다"
$ws.Range("B53").Value = "신은 존재한다"
$ws.Range("B54").Value = "플랭크를 1분이상 할 수 있다"
$ws.Range("B55").Value = "B- 재수강한다 안한다"

# Restore the view state to match where the user ended up scrolled/selected
$ws.Range("B56").Select()
